$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
# explore members
$c = $ws.Range("B6")
Write-Output ($c | Get-Member | Out-String)
